$d = $word.ActiveDocument

# The last paragraph in the document holds the _GoBack bookmark
# (bookmarkStart/bookmarkEnd) right before its paragraph mark.
# We need to insert a new paragraph containing "=2=3=13=12" just
# before that final (bookmark-only) paragraph, i.e. right after the
# text "... do one more time changes".
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Range covering just the trailing paragraph mark of the last paragraph.
$r = $lastPara.Range
[void]$r.Collapse(0)
[void]$r.MoveEnd(1, 1)

# Insert a new paragraph break, the new text, and another paragraph
# break ahead of that paragraph mark. This pushes the bookmark into a
# brand-new trailing empty paragraph, while "=2=3=13=12" becomes its
# own paragraph right before it.
$r.InsertBefore([char]13 + "=2=3=13=12" + [char]13)
